# Update the "Yearly" sheet (2017 section) with new dividend figures.
$wb = $excel.ActiveWorkbook
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# Row 5 (March, 2017): N5 taxable account value changes.
$wsYearly.Range("N5").Value = 49.69

# Row 6 (April, 2017): L6/M6/N6 values change from 0.
$wsYearly.Range("L6").Value = 44.12
$wsYearly.Range("M6").Value = 7.76
$wsYearly.Range("N6").Value = 14.81

# Recalculate so dependent formulas (O5, O6, L15:O15, All Time sheet, etc.) update.
$excel.CalculateFullRebuild()

# Update the selections / view state to match the saved workbook state.
$wsYearly.Activate()
$wsYearly.Range("O6").Select()

$wsAllTime.Activate()
$wsAllTime.Application.ActiveWindow.ScrollRow = 31
$wsAllTime.Range("K39").Select()
